# TC 51999 / 51987 / 51990 - "Add Devices" sheet test data update
#
# L10/M10 and L11/M11 were storing the text "NA" as placeholders; the test
# data now represents them as actual boolean FALSE values (the columns are
# used as checkbox-style flags for other slot cards).
$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Add Devices")

$ws.Range("L10").Value = $false
$ws.Range("M10").Value = $false
$ws.Range("L11").Value = $false
$ws.Range("M11").Value = $false

# Leave the cursor on the range that was last edited, same as the author did.
$ws.Range("L10:M11").Select()
